$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testCitizen")

$ws.Range("A2").Value = "Aleksandriya"
$ws.Range("A3").Value = "Kolonya"
$ws.Range("A4").Value = "Dollanda"
$ws.Range("A5").Value = "Fraksiya"
$ws.Range("A1").Value = "Transilvanya"
$ws.Range("A6").Value = "Bondra"
$ws.Range("A7").Value = "Kuasilinya"
$ws.Range("A8").Value = "Bulanya"

$ws.Range("B1").Select()
